$d = $word.ActiveDocument

# --- 1) Merge "G G7 G G7 G G7 " + "G " into a single run "G G7 G G7 G G7 G " ---
$r1 = $d.Content
$found1 = $r1.Find.Execute("G G7 G G7 G G7 G ", $true, $false, $false, $false, $false, $true, 1, $false, "G G7 G G7 G G7 G ", 2)

# --- 2) Merge "G7 " + " " into a single run "G7  " ---
$r2 = $d.Content
$found2 = $r2.Find.Execute("G7  ", $true, $false, $false, $false, $false, $true, 1, $false, "G7  ", 2)

# --- 3) Merge "D D7 D D7 D  D7 " + "D  D7" into a single run "D D7 D D7 D  D7 D  D7" ---
$r3 = $d.Content
$found3 = $r3.Find.Execute("D D7 D D7 D  D7 D  D7", $true, $false, $false, $false, $false, $true, 1, $false, "D D7 D D7 D  D7 D  D7", 2)

# --- 4) Move the _GoBack bookmark from the end of "Walk down:" paragraph to the
#        end of the "G7 353433" paragraph. We use the classic "type a temp char,
#        wrap a bookmark around it, then delete the char" trick so the bookmark
#        ends up as a true zero-width mark glued to the end of that paragraph's text. ---
$pG7 = $d.Paragraphs.Item(6)
$rEnd = $pG7.Range.Duplicate
[void]$rEnd.MoveEnd(1, -1)
$rEnd.Collapse(0)
$rEnd.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $rEnd)
$rEnd.Text = ""

# --- 5) Split "   CAPO on 2nd FRET or singing Dare by " into three runs:
#        "   " (bold), "CAPO on 2nd FRET" (bold + yellow highlight),
#        " or singing Dare by " (unchanged). ---
$r5 = $d.Content
$found5 = $r5.Find.Execute("   CAPO on 2nd FRET", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r5.Font.Bold = $true

$r6 = $d.Content
$r6.Find.ClearFormatting()
$r6.Find.Replacement.ClearFormatting()
$r6.Find.Replacement.Highlight = $true
$found6 = $r6.Find.Execute("CAPO on 2nd FRET", $true, $false, $false, $false, $false, $true, 1, $true, "CAPO on 2nd FRET", 2)
